$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.921.27"
$ws.Range("E2").Value = "  -0.03%  "

# Row 3
$ws.Range("D3").Value = "2.244.63"
$ws.Range("E3").Value = "  -2.25%  "

# Row 4
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").Value = "230.61"
$ws.Range("E5").Value = "  -0.58%  "

# Row 6
$ws.Range("D6").Value = "0.642"
$ws.Range("E6").Value = "  +2.01%  "

# Row 7
$ws.Range("D7").Value = "63.12"
$ws.Range("E7").Value = "  -1.87%  "

# Row 8
$ws.Range("E8").Value = "  -0.08%  "

# Row 9
$ws.Range("D9").Value = "0.447"
$ws.Range("E9").Value = "  +5.29%  "

# Row 10
$ws.Range("D10").Value = "0.0979"
$ws.Range("E10").Value = "  +1.76%  "

# Row 11
$ws.Range("D11").Value = "57.11"
$ws.Range("E11").Value = "  -1.02%  "

# Row 12
$ws.Range("D12").Value = "26.46"
$ws.Range("E12").Value = "  +0.11%  "

# Row 13
$ws.Range("D13").Value = "0.105"
$ws.Range("E13").Value = "  +0.43%  "

# Row 14
$ws.Range("D14").Value = "2.576.49"
$ws.Range("E14").Value = "  -2.24%  "

# Row 15
$ws.Range("D15").Value = "15.43"
$ws.Range("E15").Value = "  -3.58%  "

# Row 16
$ws.Range("E16").Value = "  +2.48%  "

# Row 17
$ws.Range("D17").Value = "0.825"
$ws.Range("E17").Value = "  +0.49%  "

# Row 18
$ws.Range("D18").Value = "2.242.29"
$ws.Range("E18").Value = "  -2.19%  "

# Row 19
$ws.Range("D19").Value = "43.785.04"
$ws.Range("E19").Value = "  -0.06%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0988"
$ws.Range("E20").Value = "  +2.74%  "

# Row 21
$ws.Range("D21").Value = "72.58"
$ws.Range("E21").Value = "  -1.58%  "

# Row 22
$ws.Range("D22").Value = "6.06"
$ws.Range("E22").Value = "  -2.65%  "

# Row 23
$ws.Range("D23").Value = "247.54"
$ws.Range("E23").Value = "  -3.09%  "

# Row 24
$ws.Range("E24").Value = "  +0.01%  "

# Row 25
$ws.Range("D25").Value = "2.40"
$ws.Range("E25").Value = "  -9.27%  "

# Row 26
$ws.Range("B26").Value = "WEMIXToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D26").Value = "3.34"
$ws.Range("E26").Value = "  +18.62%  "

# Row 27
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").Value = "2.24"
$ws.Range("E27").Value = "  -9.66%  "

# Row 28
$ws.Range("D28").Value = "9.82"
$ws.Range("E28").Value = "  -1.06%  "

# Row 29
$ws.Range("D29").Value = "171.66"
$ws.Range("E29").Value = "  +0.12%  "

# Row 30
$ws.Range("D30").Value = "20.93"
$ws.Range("E30").Value = "  +0.92%  "

# Row 31
$ws.Range("E31").Value = "  -1.02%  "

# Row 32
$ws.Range("E32").Value = "  -2.89%  "

# Row 33
$ws.Range("E33").Value = "  +2.16%  "

# Row 34
$ws.Range("D34").Value = "0.0687"
$ws.Range("E34").Value = "  -1.78%  "

# Row 35
$ws.Range("D35").Value = "4.77"
$ws.Range("E35").Value = "  +0.18%  "

# Row 36
$ws.Range("D36").Value = "4.94"
$ws.Range("E36").Value = "  -3.74%  "

# Row 37
$ws.Range("D37").Value = "3.64"
$ws.Range("E37").Value = "  -2.52%  "

# Row 38
$ws.Range("D38").Value = "6.42"
$ws.Range("E38").Value = "  -2.75%  "

# Row 39
$ws.Range("D39").Value = "2.27"
$ws.Range("E39").Value = "  -5.01%  "

# Row 40
$ws.Range("D40").Value = "0.0254"
$ws.Range("E40").Value = "  +1.56%  "

# Row 41
$ws.Range("E41").Value = "  -0.12%  "

# Row 42
$ws.Range("B42").Value = "TerraClassic"
$ws.Range("C42").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D42").Value = "0.000228"
$ws.Range("E42").Value = "  +2.14%  "

# Row 43
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "8.30"
$ws.Range("E43").Value = "  -5.35%  "

# Row 44
$ws.Range("D44").Value = "17.03"
$ws.Range("E44").Value = "  -1.56%  "

# Row 45
$ws.Range("D45").Value = "97.05"
$ws.Range("E45").Value = "  -2.27%  "

# Row 46
$ws.Range("E46").Value = "  -2.92%  "

# Row 47
$ws.Range("D47").Value = "0.0943"
$ws.Range("E47").Value = "  -2.77%  "

# Row 48
$ws.Range("D48").Value = "4.34"
$ws.Range("E48").Value = "  -4.18%  "

# Row 49
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "1.429.92"
$ws.Range("E49").Value = "  -4.04%  "

# Row 50
$ws.Range("B50").Value = "Celestia"
$ws.Range("C50").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D50").Value = "9.86"
$ws.Range("E50").Value = "  -10.53%  "

# Row 51
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "2.29"
$ws.Range("E51").Value = "  -1.01%  "
